$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: feature number bumps from 9 to 10, and its "Known Issues / Notes"
# placeholder text ("7`n") is replaced with the real note.
$ws.Range("B11").Value = 10
$ws.Range("C11").Value = "User Entry- Init of New Entry"

# New row 12: the item that used to be numbered 9 now documents the new
# "Task entry form" feature.
$ws.Range("B12").Value = 9
$ws.Range("C12").Value = "Task entry form"

# The manifest-entry description in D9 no longer calls bClearUserEntryNew();
# the method call is truncated since "Hour rounds off now" replaces it.
$ws.Range("D9").Value = "manifest entry - android:name=`".CustomClasses.taskQGlobal`"`n`n((taskQGlobal) getApplication())."

# Move the active selection to H9, matching the saved view state.
$ws.Range("H9").Select() | Out-Null
